$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the password cell's value (B3) from "admin" to "manager"
$ws.Range("B3").Value = "manager"

# Move/update the current selection to B3 (matches the saved view state)
$ws.Range("B3").Select()
